$d = $word.ActiveDocument

$replacements = @(
    @("153×3=", "263×5="),
    @("701×5=", "461×2="),
    @("499×4=", "767×5="),
    @("120×9=", "751×5="),
    @("534×9=", "794×4="),
    @("493×2=", "858×6="),
    @("271×5=", "487×9="),
    @("942×3=", "562×6="),
    @("659×4=", "804×8="),
    @("424×8=", "415×4="),
    @("502×3=", "304×7="),
    @("408×8=", "836×5="),
    @("510×8=", "460×4="),
    @("785×5=", "857×9="),
    @("127×8=", "651×9="),
    @("310×9=", "240×4="),
    @("608×7=", "649×9="),
    @("188×7=", "390×4="),
    @("600×4=", "112×4="),
    @("930×6=", "460×2="),
    @("878×9=", "538×6="),
    @("524×6=", "831×9="),
    @("584×2=", "176×3="),
    @("648×3=", "476×9="),
    @("261×7=", "751×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
